# Applies the "TaskManager vorerst abgeschlossen, mit Qt Creator Projekt
# erstellt" edit to the architecture document:
#   1. Adds a new "Qt Designer zum erstellen der grafischen Oberfläche"
#      bullet to the "Entwicklungswerkzeuge" list, right before the
#      "CMake als Buildsystem" bullet.
#   2. Removes the still-unimplemented "FileHandler" stub bullets, and
#      the (otherwise blank) spacer paragraphs that used to sit right
#      above the "Struktur" and "3.2 Abhängigkeiten" headings.

$d = $word.ActiveDocument

function Get-ParaText($para) {
    # Paragraph.Range.Text includes the trailing paragraph-mark / cell-mark
    # character(s); strip those so comparisons are exact.
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ((Get-ParaText $d.Paragraphs($i)) -eq $text) {
            return $i
        }
    }
    return -1
}

function Remove-ParaWithText($text) {
    $idx = Find-ParaIndex $text
    if ($idx -gt 0) {
        $d.Paragraphs($idx).Range.Delete()
    }
}

function Remove-BlankParaBefore($text) {
    $idx = Find-ParaIndex $text
    if ($idx -gt 1) {
        $above = $d.Paragraphs($idx - 1)
        if ((Get-ParaText $above) -eq "") {
            $above.Range.Delete()
        }
    }
}

# --- 1. Insert the new "Qt Designer" bullet before "CMake als Buildsystem" ---

$cmakeIdx = Find-ParaIndex "CMake als Buildsystem"
$cmakePara = $d.Paragraphs($cmakeIdx)
$cmakePara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($cmakeIdx)
$newPara.Range.Text = "Qt Designer zum erstellen der grafischen Oberfläche"

# --- 2. Drop the spacer paragraph above "Struktur" ---

Remove-BlankParaBefore "Struktur"

# --- 3. Remove the "FileHandler" stub bullet (3.1 Hauptkomponenten) ---

Remove-ParaWithText "FileHandler"

# --- 4. Drop the spacer paragraph above "3.2 Abhängigkeiten" ---

Remove-BlankParaBefore "3.2 Abhängigkeiten"

# --- 5. Remove the "FileHandler wird von TaskManager genutzt" stub bullet ---

Remove-ParaWithText "FileHandler wird von TaskManager genutzt"
